$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Material")

# 1) Update row 4 (Trallgolv 120mm -> Trallgolv, lm -> m2, new quantities)
$ws.Range("B4").Value = "Trallgolv"
$ws.Range("C4").Value = "m2"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.15
$ws.Range("F4").Value = 195
$ws.Range("G4").Value = 350
$ws.Range("I4").Value = "Beräknas på golvyta"

# 2) Insert a new row at position 6 (before "Spånskivegolv") for "Trossbottenpapp"
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "Golv"
$ws.Range("B6").Value = "Trossbottenpapp"
$ws.Range("C6").Value = "m2"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.05
$ws.Range("F6").Value = 25
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = $true
$ws.Range("I6").Value = "Beräknas på golvyta"

# 3) The rows that shifted down (old 6 "Spånskivegolv" -> new 7, old 7 "Isolering golv" -> new 8)
#    now also get a note added.
$ws.Range("I7").Value = "Beräknas på golvyta"
$ws.Range("I8").Value = "Beräknas på golvyta"
